$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 244
$ws1.Range("F5").Value = 2869
$ws1.Range("F6").Value = 1984
$ws1.Range("F7").Value = 382
$ws1.Range("F9").Value = 1054
$ws1.Range("F10").Value = 198
$ws1.Range("F11").Value = 336
$ws1.Range("F12").Value = 47

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 244
$ws4.Range("F5").Value = 2869
$ws4.Range("F6").Value = 1984
$ws4.Range("F7").Value = 382
$ws4.Range("F10").Value = 1054
$ws4.Range("F11").Value = 198
$ws4.Range("F12").Value = 336
$ws4.Range("F13").Value = 47
